$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (C) column for all existing data rows (2..479)
#    from 45175 (2023-09-06) to 45177 (2023-09-08).
for ($r = 2; $r -le 479; $r++) {
    $ws.Cells.Item($r, 3).Value = 45177
}

# 2) Row 479 gains an explicit row height (matches every other data row).
$ws.Rows.Item(479).RowHeight = 15

# 3) Append the four new cases reported for 2023-09-07 / changed 2023-09-08.
$newRows = @(
    @{ Row = 480; A = "A 41724-2023"; G = 2.1 },
    @{ Row = 481; A = "A 41739-2023"; G = 1.3 },
    @{ Row = 482; A = "A 41744-2023"; G = 2.4 },
    @{ Row = 483; A = "A 41746-2023"; G = 1.5 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = 45176
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 3).Value = 45177
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 4).Value = "JÖNKÖPINGS LÄN"
    $ws.Cells.Item($r, 5).Value = "EKSJÖ"
    $ws.Cells.Item($r, 7).Value = $entry.G
    for ($c = 8; $c -le 17; $c++) {
        $ws.Cells.Item($r, $c).Value = 0
    }
    $ws.Cells.Item($r, 18).Value = ""
    $ws.Cells.Item($r, 18).WrapText = $true
}

# Rows 480-482 keep the standard explicit row height; row 483 (the new
# last row) is left without one, matching how row 479 looked before this
# edit (the final row of the sheet has no explicit height).
$ws.Rows.Item(480).RowHeight = 15
$ws.Rows.Item(481).RowHeight = 15
$ws.Rows.Item(482).RowHeight = 15
